# Edit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Replace the (single) previous worker's account-statement detail rows with a
# new, larger set of workers/periods, and update the summary header figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the summary block above the table (RAZON SOCIAL / NIT unchanged,
#    VALOR MORA total and worker/period counters change).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1070416      # VALOR MORA
$ws.Range("C13").Value = 5            # Cant. Trabajadores
$ws.Range("F13").Value = 8            # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) Grow the detail table from 2 rows (16:17) to 26 rows (16:41), preserving
#    the existing look: row 16 supplies the "interior" row style, the old
#    row 17 supplies the "bottom border / last row" style which must end up
#    on the new last data row (41). Inserting whole rows before the old
#    last row pushes it down and shifts the merged cells / signature block
#    below (previously 22:23) down to 46:47 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("17:40").Insert()

$templateRow = $ws.Range("B16:J16")
$fillRange = $ws.Range("B17:J40")
$templateRow.Copy($fillRange)

# ---------------------------------------------------------------------------
# 3) Write the new worker / period / value data into rows 16-41.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------------
$rowsData = @(
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2507", 35112, 877803),
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2506", 35112, 877803),
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2505", 35112, 877803),
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2504", 35112, 877803),
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2503", 35112, 877803),
    @("CC", "1047489453", "YAIR JOSE GUERRERO CARABALLO", "2502", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2507", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2506", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2505", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2504", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2503", 35112, 877803),
    @("CC", "1047510253", "JEAN PAUL NARVAEZ HERNANDEZ", "2502", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2507", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2506", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2505", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2504", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2503", 35112, 877803),
    @("CC", "1235039234", "LISSET LORENA ATENCIO BARRAGAN", "2502", 35112, 877803),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2507", 46400, 1160000),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2506", 46400, 1160000),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2505", 46400, 1160000),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2504", 46400, 1160000),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2503", 46400, 1160000),
    @("CC", "1007901612", "SEBASTIAN MARTINEZ MONTES", "2502", 46400, 1160000),
    @("CC", "1007763932", "SOFIA CATALINA ALANDETE GONZALEZ", "2204", 80000, 2000000),
    @("CC", "1007763932", "SOFIA CATALINA ALANDETE GONZALEZ", "2203", 80000, 2000000)
)

$startRow = 16
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $row = $rowsData[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]   # B Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[1]   # C N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G Salario Basico
}
